$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-9 contain numeric Cash/Credit/Cashless (mobile) figures in columns B, C, D
# that were entered as cents-like integers (one decimal place too large). Rescale
# them down by a factor of 10, then recompute the dependent Total (E) and
# Payable (F) columns from the corrected figures.
for ($row = 2; $row -le 9; $row++) {
    $b = $ws.Cells.Item($row, 2).Value2
    $c = $ws.Cells.Item($row, 3).Value2
    $d = $ws.Cells.Item($row, 4).Value2

    $newB = $b / 10
    $newC = $c / 10
    $newD = $d / 10

    $ws.Cells.Item($row, 2).Value2 = $newB
    $ws.Cells.Item($row, 3).Value2 = $newC
    $ws.Cells.Item($row, 4).Value2 = $newD

    $newE = $newB + $newC + $newD
    $ws.Cells.Item($row, 5).Value2 = $newE
    $ws.Cells.Item($row, 6).Value2 = [Math]::Round($newE / 2)
}
